$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1325.3334
$ws.Range("J32").Value = 1275
$ws.Range("L32").Value = 1275
$ws.Range("N32").Value = -1927
$ws.Range("H53").Value = 224.875
$ws.Range("I53").Value = 132.66667
$ws.Range("K53").Value = 132.66667
$ws.Range("M53").Value = 504.33333
$ws.Range("H87").Value = 74086.45
$ws.Range("J87").Value = 88868.875
$ws.Range("L87").Value = 88868.875
$ws.Range("N87").Value = -91364.875
$ws.Range("H90").Value = 74086.45
$ws.Range("J90").Value = 88868.875
$ws.Range("L90").Value = 266606.625
$ws.Range("N90").Value = -279086.625
$ws.Range("H98").Value = 2346.5881
$ws.Range("I98").Value = 2570
$ws.Range("J98").Value = 2027.4286
$ws.Range("K98").Value = 2570
$ws.Range("L98").Value = 2027.4286
$ws.Range("M98").Value = -1072
$ws.Range("N98").Value = -5023.4286
$ws.Range("H122").Value = 2346.5881
$ws.Range("I122").Value = 2570
$ws.Range("J122").Value = 2027.4286
$ws.Range("K122").Value = 7710
$ws.Range("L122").Value = 6082.2858
$ws.Range("M122").Value = -5260
$ws.Range("N122").Value = -10982.2858
$ws.Range("H128").Value = 80000
$ws.Range("J128").Value = 80000
$ws.Range("L128").Value = 80000
$ws.Range("N128").Value = -89960
$ws.Range("H129").Value = 2287.1333
$ws.Range("I129").Value = 1525.6666
$ws.Range("K129").Value = 4576.9998
$ws.Range("M129").Value = 423.0002000000004
$ws.Range("H132").Value = 1709.125
$ws.Range("I132").Value = 1674.7609
$ws.Range("K132").Value = 5024.2827
$ws.Range("M132").Value = -2494.2827
$ws.Range("H137").Value = 2447.5925
$ws.Range("I137").Value = 1248.6666
$ws.Range("J137").Value = 2790.1428
$ws.Range("K137").Value = 3745.9998
$ws.Range("L137").Value = 8370.428400000001
$ws.Range("M137").Value = -1195.9998
$ws.Range("N137").Value = -13470.4284
$ws.Range("H138").Value = 2575.6667
$ws.Range("I138").Value = 1698.4445
$ws.Range("J138").Value = 3452.889
$ws.Range("K138").Value = 5095.333500000001
$ws.Range("L138").Value = 10358.667
$ws.Range("M138").Value = 44.66649999999936
$ws.Range("N138").Value = -20638.667

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3434.3396
$ws.Range("I132").Value = 2469.6956
$ws.Range("K132").Value = 7409.0868
$ws.Range("M132").Value = -4879.0868
$ws.Range("H139").Value = 84999.836
$ws.Range("J139").Value = 84999.836
$ws.Range("L139").Value = 84999.836
$ws.Range("N139").Value = -95279.836

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1903.35
$ws.Range("I20").Value = 2143.889
$ws.Range("J20").Value = 1706.5454
$ws.Range("K20").Value = 2143.889
$ws.Range("L20").Value = 1706.5454
$ws.Range("M20").Value = -1896.889
$ws.Range("N20").Value = -2200.5454
$ws.Range("H86").Value = 1172.2354
$ws.Range("I86").Value = 1078.4546
$ws.Range("J86").Value = 1344.1666
$ws.Range("K86").Value = 1078.4546
$ws.Range("L86").Value = 1344.1666
$ws.Range("M86").Value = 44.54539999999997
$ws.Range("N86").Value = -3590.1666
$ws.Range("H89").Value = 1172.2354
$ws.Range("I89").Value = 1078.4546
$ws.Range("J89").Value = 1344.1666
$ws.Range("K89").Value = 5392.273
$ws.Range("L89").Value = 6720.833000000001
$ws.Range("M89").Value = 223.7269999999999
$ws.Range("N89").Value = -17952.833
$ws.Range("H134").Value = 3987.6904
$ws.Range("I134").Value = 2071.5
$ws.Range("K134").Value = 6214.5
$ws.Range("M134").Value = -3679.5
$ws.Range("H140").Value = 174111.17
$ws.Range("J140").Value = 174111.17
$ws.Range("L140").Value = 174111.17
$ws.Range("N140").Value = -184471.17

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 781.0909
$ws.Range("I16").Value = 659.2
$ws.Range("K16").Value = 659.2
$ws.Range("M16").Value = -372.2
$ws.Range("H31").Value = 5173.965
$ws.Range("I31").Value = 2718.2144
$ws.Range("J31").Value = 7545.0347
$ws.Range("K31").Value = 2718.2144
$ws.Range("L31").Value = 7545.0347
$ws.Range("M31").Value = -2423.2144
$ws.Range("N31").Value = -8135.0347
$ws.Range("H34").Value = 5173.965
$ws.Range("I34").Value = 2718.2144
$ws.Range("J34").Value = 7545.0347
$ws.Range("K34").Value = 2718.2144
$ws.Range("L34").Value = 7545.0347
$ws.Range("M34").Value = -2516.2144
$ws.Range("N34").Value = -7949.0347
$ws.Range("H68").Value = 66980.28999999999
$ws.Range("J68").Value = 69977
$ws.Range("L68").Value = 69977
$ws.Range("N68").Value = -71475
$ws.Range("H71").Value = 66980.28999999999
$ws.Range("J71").Value = 69977
$ws.Range("L71").Value = 209931
$ws.Range("N71").Value = -217419
$ws.Range("H107").Value = 984.75
$ws.Range("I107").Value = 993.94116
$ws.Range("J107").Value = 932.6667
$ws.Range("K107").Value = 993.94116
$ws.Range("L107").Value = 932.6667
$ws.Range("M107").Value = 926.05884
$ws.Range("N107").Value = -4772.6667
$ws.Range("H113").Value = 781.0909
$ws.Range("I113").Value = 659.2
$ws.Range("K113").Value = 659.2
$ws.Range("M113").Value = 1510.8
$ws.Range("H122").Value = 3243.6875
$ws.Range("J122").Value = 4104.8184
$ws.Range("L122").Value = 12314.4552
$ws.Range("N122").Value = -17214.4552
$ws.Range("H134").Value = 1731.1578
$ws.Range("J134").Value = 1949
$ws.Range("L134").Value = 5847
$ws.Range("N134").Value = -10917

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 11827.571
$ws.Range("I87").Value = 5951.0713
$ws.Range("K87").Value = 17853.2139
$ws.Range("M87").Value = -16605.2139
$ws.Range("H90").Value = 11827.571
$ws.Range("I90").Value = 5951.0713
$ws.Range("K90").Value = 53559.64169999999
$ws.Range("M90").Value = -47319.64169999999
$ws.Range("H109").Value = 3946.5293
$ws.Range("I109").Value = 1926.875
$ws.Range("J109").Value = 5741.778
$ws.Range("K109").Value = 5780.625
$ws.Range("L109").Value = 17225.334
$ws.Range("M109").Value = -4740.625
$ws.Range("N109").Value = -19305.334

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 43750
$ws.Range("J26").Value = 43750
$ws.Range("L26").Value = 43750
$ws.Range("N26").Value = -44310
$ws.Range("H47").Value = 29465.5
$ws.Range("J47").Value = 29465.5
$ws.Range("L47").Value = 29465.5
$ws.Range("N47").Value = -30601.5
$ws.Range("H50").Value = 43750
$ws.Range("J50").Value = 43750
$ws.Range("L50").Value = 43750
$ws.Range("N50").Value = -44746
$ws.Range("H70").Value = 27998.6
$ws.Range("J70").Value = 9998
$ws.Range("L70").Value = 9998
$ws.Range("N70").Value = -10538
$ws.Range("H73").Value = 27998.6
$ws.Range("J73").Value = 9998
$ws.Range("L73").Value = 9998
$ws.Range("N73").Value = -11870

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 145416.58
$ws.Range("I22").Value = 335973.66
$ws.Range("K22").Value = 335973.66
$ws.Range("M22").Value = -335678.66
$ws.Range("H27").Value = 145416.58
$ws.Range("I27").Value = 335973.66
$ws.Range("K27").Value = 335973.66
$ws.Range("M27").Value = -335866.66
$ws.Range("H46").Value = 2780.2903
$ws.Range("J46").Value = 3444.2273
$ws.Range("L46").Value = 3444.2273
$ws.Range("N46").Value = -3820.2273
$ws.Range("H62").Value = 63963.168
$ws.Range("J62").Value = 71998
$ws.Range("L62").Value = 71998
$ws.Range("N62").Value = -73246
$ws.Range("H65").Value = 63963.168
$ws.Range("J65").Value = 71998
$ws.Range("L65").Value = 215994
$ws.Range("N65").Value = -222234

Write-Output "Applied all cell updates."